$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 100, shifting the
# existing rows 100-194 down to 102-196.
$ws.Rows.Item(100).Resize(2).Insert()

# Populate the two newly inserted rows (100 and 101) with fresh data.
# Row 100
$ws.Cells.Item(100, 1).Value2 = 11
$ws.Cells.Item(100, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(100, 3).Value2 = "Bíobío"
$ws.Cells.Item(100, 4).Value2 = 44586
$ws.Cells.Item(100, 5).Value2 = 8
$ws.Cells.Item(100, 6).Value2 = 100112017
$ws.Cells.Item(100, 7).Value2 = "Apio"
$ws.Cells.Item(100, 8).Value2 = "Americana (o)"
$ws.Cells.Item(100, 9).Value2 = "Primera"
$ws.Cells.Item(100, 10).Value2 = 280
$ws.Cells.Item(100, 11).Value2 = 6500
$ws.Cells.Item(100, 12).Value2 = 7000
$ws.Cells.Item(100, 13).Value2 = 6732
$ws.Cells.Item(100, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(100, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(100, 16).Value2 = 1122
$ws.Cells.Item(100, 17).Value2 = 6
$ws.Cells.Item(100, 18).Value2 = "Hortaliza"

# Row 101
$ws.Cells.Item(101, 1).Value2 = 11
$ws.Cells.Item(101, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(101, 3).Value2 = "Bíobío"
$ws.Cells.Item(101, 4).Value2 = 44586
$ws.Cells.Item(101, 5).Value2 = 8
$ws.Cells.Item(101, 6).Value2 = 100112017
$ws.Cells.Item(101, 7).Value2 = "Apio"
$ws.Cells.Item(101, 8).Value2 = "Americana (o)"
$ws.Cells.Item(101, 9).Value2 = "Segunda"
$ws.Cells.Item(101, 10).Value2 = 190
$ws.Cells.Item(101, 11).Value2 = 4500
$ws.Cells.Item(101, 12).Value2 = 5000
$ws.Cells.Item(101, 13).Value2 = 4737
$ws.Cells.Item(101, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(101, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(101, 16).Value2 = 790
$ws.Cells.Item(101, 17).Value2 = 6
$ws.Cells.Item(101, 18).Value2 = "Hortaliza"

# Match the number format of column D used throughout the table for the
# date cells in the two new rows.
$ws.Cells.Item(100, 4).NumberFormat = $ws.Cells.Item(99, 4).NumberFormat
$ws.Cells.Item(101, 4).NumberFormat = $ws.Cells.Item(99, 4).NumberFormat
